$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header label for the row-name column
$ws.Range("A1").Value = "Variant"

# Rename first and last mutant rows
$ws.Range("A2").Value = "Wild type"
$ws.Range("A9").Value = "Triple"

# Update selection to match the committed state
$ws.Range("B5").Select()
